# Update NATMI ligand-receptor TPM-derived values on Sheet1
# (recomputed specificity/weight metrics for rows 2-9)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.04433703455491324
$ws.Range("J2").Value = 0.04433703455491323
$ws.Range("M2").Value = 0.3754836666666666
$ws.Range("N2").Value = 1.126451
$ws.Range("O2").Value = 0.1282941949732911
$ws.Range("P2").Value = 0.1282941949732911
$ws.Range("Q2").Value = 0.03452910250299999
$ws.Range("R2").Value = 0.310761922527
$ws.Range("S2").Value = 0.005688184155725583
$ws.Range("T2").Value = 0.005688184155725582
$ws.Range("I3").Value = 0.04433703455491324
$ws.Range("J3").Value = 0.04433703455491323
$ws.Range("O3").Value = 0.1860305746394907
$ws.Range("P3").Value = 0.1860305746394907
$ws.Range("S3").Value = 0.008248044016061468
$ws.Range("T3").Value = 0.008248044016061466
$ws.Range("I4").Value = 0.04433703455491324
$ws.Range("J4").Value = 0.04433703455491323
$ws.Range("M4").Value = 0.9905483333333333
$ws.Range("N4").Value = 2.971645
$ws.Range("O4").Value = 0.3384477469693805
$ws.Range("P4").Value = 0.3384477469693805
$ws.Range("Q4").Value = 0.091089834185
$ws.Range("R4").Value = 0.8198085076649999
$ws.Range("S4").Value = 0.01500576945241395
$ws.Range("T4").Value = 0.01500576945241395
$ws.Range("I5").Value = 0.04433703455491324
$ws.Range("J5").Value = 0.04433703455491323
$ws.Range("M5").Value = 1.016244333333333
$ws.Range("N5").Value = 3.048733
$ws.Range("O5").Value = 0.3472274834178377
$ws.Range("P5").Value = 0.3472274834178377
$ws.Range("Q5").Value = 0.093452812649
$ws.Range("R5").Value = 0.8410753138409999
$ws.Range("S5").Value = 0.01539503693071223
$ws.Range("T5").Value = 0.01539503693071223
$ws.Range("G6").Value = 1.982131
$ws.Range("H6").Value = 5.946393
$ws.Range("I6").Value = 0.9556629654450868
$ws.Range("J6").Value = 0.9556629654450867
$ws.Range("M6").Value = 0.3754836666666666
$ws.Range("N6").Value = 1.126451
$ws.Range("O6").Value = 0.1282941949732911
$ws.Range("P6").Value = 0.1282941949732911
$ws.Range("Q6").Value = 0.7442578156936666
$ws.Range("R6").Value = 6.698320341243
$ws.Range("S6").Value = 0.1226060108175655
$ws.Range("T6").Value = 0.1226060108175655
$ws.Range("G7").Value = 1.982131
$ws.Range("H7").Value = 5.946393
$ws.Range("I7").Value = 0.9556629654450868
$ws.Range("J7").Value = 0.9556629654450867
$ws.Range("O7").Value = 0.1860305746394907
$ws.Range("P7").Value = 0.1860305746394907
$ws.Range("Q7").Value = 1.079196990653
$ws.Range("R7").Value = 9.712772915877
$ws.Range("S7").Value = 0.1777825306234293
$ws.Range("T7").Value = 0.1777825306234292
$ws.Range("G8").Value = 1.982131
$ws.Range("H8").Value = 5.946393
$ws.Range("I8").Value = 0.9556629654450868
$ws.Range("J8").Value = 0.9556629654450867
$ws.Range("M8").Value = 0.9905483333333333
$ws.Range("N8").Value = 2.971645
$ws.Range("O8").Value = 0.3384477469693805
$ws.Range("P8").Value = 0.3384477469693805
$ws.Range("Q8").Value = 1.963396558498333
$ws.Range("R8").Value = 17.670569026485
$ws.Range("S8").Value = 0.3234419775169665
$ws.Range("T8").Value = 0.3234419775169665
$ws.Range("G9").Value = 1.982131
$ws.Range("H9").Value = 5.946393
$ws.Range("I9").Value = 0.9556629654450868
$ws.Range("J9").Value = 0.9556629654450867
$ws.Range("M9").Value = 1.016244333333333
$ws.Range("N9").Value = 3.048733
$ws.Range("O9").Value = 0.3472274834178377
$ws.Range("P9").Value = 0.3472274834178377
$ws.Range("Q9").Value = 2.014329396674333
$ws.Range("R9").Value = 18.128964570069
$ws.Range("S9").Value = 0.3318324464871255
$ws.Range("T9").Value = 0.3318324464871255
